$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1480.2
$ws.Range("J40").Value = 1466.6666
$ws.Range("L40").Value = 1466.6666
$ws.Range("N40").Value = -1816.6666
$ws.Range("H76").Value = 3327.2727
$ws.Range("I76").Value = 3325
$ws.Range("J76").Value = 3333.3333
$ws.Range("K76").Value = 3325
$ws.Range("L76").Value = 3333.3333
$ws.Range("M76").Value = -3010
$ws.Range("N76").Value = -3963.3333
$ws.Range("H79").Value = 3327.2727
$ws.Range("I79").Value = 3325
$ws.Range("J79").Value = 3333.3333
$ws.Range("K79").Value = 3325
$ws.Range("L79").Value = 3333.3333
$ws.Range("M79").Value = -2233
$ws.Range("N79").Value = -5517.3333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 12625.75
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 12625.75
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -19133.75
$ws.Range("H138").Value = 2332.192
$ws.Range("I138").Value = 881.5599999999999
$ws.Range("J138").Value = 2822.2703
$ws.Range("K138").Value = 2644.68
$ws.Range("L138").Value = 8466.8109
$ws.Range("M138").Value = 2495.32
$ws.Range("N138").Value = -18746.8109

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5755.0366
$ws.Range("I32").Value = 3463.5322
$ws.Range("J32").Value = 12858.7
$ws.Range("K32").Value = 3463.5322
$ws.Range("L32").Value = 12858.7
$ws.Range("M32").Value = -3176.5322
$ws.Range("N32").Value = -13432.7
$ws.Range("H133").Value = 39286
$ws.Range("J133").Value = 39286
$ws.Range("L133").Value = 39286
$ws.Range("N133").Value = -44346
$ws.Range("H135").Value = 44500
$ws.Range("J135").Value = 44500
$ws.Range("L135").Value = 44500
$ws.Range("N135").Value = -54640

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 10980.667
$ws.Range("I26").Value = 10980.667
$ws.Range("K26").Value = 10980.667
$ws.Range("M26").Value = -10688.667
$ws.Range("H134").Value = 3200.524
$ws.Range("I134").Value = 1908.8462
$ws.Range("K134").Value = 5726.5386
$ws.Range("M134").Value = -3191.5386

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10206407
$ws.Range("I31").Value = 974.4054
$ws.Range("J31").Value = 41673160
$ws.Range("K31").Value = 974.4054
$ws.Range("L31").Value = 41673160
$ws.Range("M31").Value = -679.4054
$ws.Range("N31").Value = -41673750
$ws.Range("H34").Value = 10206407
$ws.Range("I34").Value = 974.4054
$ws.Range("J34").Value = 41673160
$ws.Range("K34").Value = 974.4054
$ws.Range("L34").Value = 41673160
$ws.Range("M34").Value = -772.4054
$ws.Range("N34").Value = -41673564
$ws.Range("H41").Value = 33518.145
$ws.Range("J41").Value = 38428
$ws.Range("L41").Value = 38428
$ws.Range("N41").Value = -39284
$ws.Range("H50").Value = 26037.223
$ws.Range("J50").Value = 26037.223
$ws.Range("L50").Value = 26037.223
$ws.Range("N50").Value = -27287.223
$ws.Range("H59").Value = 29923.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 29923.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 29923.5
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -32213.5
$ws.Range("H68").Value = 50150.777
$ws.Range("J68").Value = 50150.777
$ws.Range("L68").Value = 50150.777
$ws.Range("N68").Value = -51648.777
$ws.Range("H71").Value = 50150.777
$ws.Range("J71").Value = 50150.777
$ws.Range("L71").Value = 150452.331
$ws.Range("N71").Value = -157940.331
$ws.Range("H74").Value = 34509.75
$ws.Range("J74").Value = 34509.75
$ws.Range("L74").Value = 34509.75
$ws.Range("N74").Value = -36257.75
$ws.Range("H77").Value = 34509.75
$ws.Range("J77").Value = 34509.75
$ws.Range("L77").Value = 103529.25
$ws.Range("N77").Value = -112265.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1356.7407
$ws.Range("I5").Value = 356
$ws.Range("K5").Value = 1068
$ws.Range("M5").Value = -956
$ws.Range("H122").Value = 3406.36
$ws.Range("I122").Value = 577.73334
$ws.Range("J122").Value = 3905.5293
$ws.Range("K122").Value = 5199.60006
$ws.Range("L122").Value = 35149.7637
$ws.Range("M122").Value = -2749.60006
$ws.Range("N122").Value = -40049.7637
$ws.Range("H135").Value = 1356.7407
$ws.Range("I135").Value = 356
$ws.Range("K135").Value = 3204
$ws.Range("M135").Value = -669

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 35716884
$ws.Range("I80").Value = 62502320
$ws.Range("J80").Value = 2966.6667
$ws.Range("K80").Value = 62502320
$ws.Range("L80").Value = 2966.6667
$ws.Range("M80").Value = -62501322
$ws.Range("N80").Value = -4962.6667
$ws.Range("H83").Value = 35716884
$ws.Range("I83").Value = 62502320
$ws.Range("J83").Value = 2966.6667
$ws.Range("K83").Value = 312511600
$ws.Range("L83").Value = 14833.3335
$ws.Range("M83").Value = -312506608
$ws.Range("N83").Value = -24817.3335
$ws.Range("H97").Value = 709.1111
$ws.Range("I97").Value = 633.17645
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 633.17645
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -137.17645
$ws.Range("N97").Value = -2992
$ws.Range("H122").Value = 3283.9167
$ws.Range("I122").Value = 2163.375
$ws.Range("K122").Value = 6490.125
$ws.Range("M122").Value = -4040.125
$ws.Range("H123").Value = 10198.97
$ws.Range("J123").Value = 10198.97
$ws.Range("L123").Value = 10198.97
$ws.Range("N123").Value = -15098.97
$ws.Range("H132").Value = 5549.625
$ws.Range("I132").Value = 3100
$ws.Range("J132").Value = 6366.1665
$ws.Range("K132").Value = 9300
$ws.Range("L132").Value = 19098.4995
$ws.Range("M132").Value = -6770
$ws.Range("N132").Value = -24158.4995
$ws.Range("H134").Value = 46733.16
$ws.Range("J134").Value = 47876.375
$ws.Range("L134").Value = 143629.125
$ws.Range("N134").Value = -148699.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7577178.5
$ws.Range("I22").Value = 14706737
$ws.Range("J22").Value = 2022.5
$ws.Range("K22").Value = 14706737
$ws.Range("L22").Value = 2022.5
$ws.Range("M22").Value = -14706442
$ws.Range("N22").Value = -2612.5
$ws.Range("H27").Value = 7577178.5
$ws.Range("I27").Value = 14706737
$ws.Range("J27").Value = 2022.5
$ws.Range("K27").Value = 14706737
$ws.Range("L27").Value = 2022.5
$ws.Range("M27").Value = -14706630
$ws.Range("N27").Value = -2236.5
$ws.Range("H46").Value = 1918.8462
$ws.Range("I46").Value = 1588.3334
$ws.Range("J46").Value = 2662.5
$ws.Range("K46").Value = 1588.3334
$ws.Range("L46").Value = 2662.5
$ws.Range("M46").Value = -1400.3334
$ws.Range("N46").Value = -3038.5
$ws.Range("H61").Value = 2069.7144
$ws.Range("I61").Value = 1697.6
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1697.6
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1495.6
$ws.Range("N61").Value = -3404
$ws.Range("H113").Value = 2069.7144
$ws.Range("I113").Value = 1697.6
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1697.6
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 472.4000000000001
$ws.Range("N113").Value = -7340
$ws.Range("H122").Value = 6723.75
$ws.Range("I122").Value = 3758
$ws.Range("K122").Value = 11274
$ws.Range("M122").Value = -8824

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3989.1853
$ws.Range("I122").Value = 2747.5293
$ws.Range("K122").Value = 8242.5879
$ws.Range("M122").Value = -5792.5879
